$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A26").Value = 0.84730412960844359
$wsAbout.Range("B26").Value = "2021 dollars per 2012 dollar"
$wsAbout.Range("A18").Value = "billion 2021 dollars"
$wsAbout.Range("A21").Value = "million 2021 dollars"
$wsAbout.Range("A24").Value = "2021 dollars"
$wsAbout.Range("B29").Value = 'which in this case is "2012 dollars per 2021 dollar."'

$wsAbout.Range("A9").ClearFormats() | Out-Null

$wsSOCU = $wb.Worksheets.Item("OCCF-DpSOCU")
$wsSOCU.Range("B2").ClearFormats() | Out-Null

$wsAbout.Activate() | Out-Null
$wsAbout.Range("B30").Select() | Out-Null
